# Updated cryptos list - refresh Price and Volume(1h) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.553.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "'1.908.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "'325.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'0.4848"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.87%  "
$ws.Range("D8").Value = "'0.4079"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.08142"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("D10").Value = "'1.011"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").Value = "'23.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.32%  "
$ws.Range("D12").Value = "'6.032"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "'1.854.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").Value = "'7.101"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "'90.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "'0.00001044"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "'17.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'29.563.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").Value = "'11.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "'2.125.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "'154.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").Value = "'20.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("D28").Value = "'6.330"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.13%  "
$ws.Range("D29").Value = "'2.106"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").Value = "'119.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").Value = "'1.038"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("D32").Value = "'0.09566"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "'5.546"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").Value = "'1.398"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("D35").Value = "'3.556"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").Value = "'0.02268"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'0.06120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'1.173"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'0.5966"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("D40").Value = "'7.954"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.60%  "
$ws.Range("D41").Value = "'10.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.94%  "
$ws.Range("D42").Value = "'0.1859"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "'2.443"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "'1.280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").Value = "'0.07728"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "'12.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").Value = "'1.956"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").Value = "'115.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'72.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("E51").Value = "  +2.64%  "
